# Update cryptos list price/volume figures (GitHub Actions data refresh).
# For numeric-looking Price values we force the cell to Text ("@") before
# assigning, then reset the style back to "Normal" afterwards - this keeps
# the literal text (e.g. "580.04", "1.00", "42.30") exactly as scraped
# instead of letting Excel's auto-type-detection coerce it into a binary
# double (which would both lose trailing zeros and introduce float noise).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.526.01'
$ws.Range("E2").Value = '  +0.84%  '
$ws.Range("D3").Value = '3.445.69'
$ws.Range("E3").Value = '  +1.65%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.04'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.74'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +8.89%  '
$ws.Range("D7").Value = '3.447.36'
$ws.Range("E7").Value = '  +1.69%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.474'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.75%  '
$ws.Range("E10").Value = '  +2.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.126'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.45%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.391'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.05%  '
$ws.Range("D13").Value = '4.033.13'
$ws.Range("E13").Value = '  +1.65%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.93'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.72%  '
$ws.Range("E15").Value = '  -0.42%  '
$ws.Range("E16").Value = '  +0.41%  '
$ws.Range("D17").Value = '3.448.60'
$ws.Range("E17").Value = '  +1.73%  '
$ws.Range("D18").Value = '61.642.42'
$ws.Range("E18").Value = '  +0.81%  '
$ws.Range("E19").Value = '  +8.58%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.36'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.49'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.20%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '389.02'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.24%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.569'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.43%  '
$ws.Range("D24").Value = '3.585.10'
$ws.Range("E24").Value = '  +1.64%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.98'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.53%  '
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.36%  '
$ws.Range("B27").Value = 'LEO'
$ws.Range("C27").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.78'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.63%  '
$ws.Range("E28").Value = '  -1.44%  '
$ws.Range("E29").Value = '  +4.21%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.76'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.90%  '
$ws.Range("E31").Value = '  +0.07%  '
$ws.Range("E32").Value = '  -12.96%  '
$ws.Range("E33").Value = '  +1.00%  '
$ws.Range("E34").Value = '  +0.89%  '
$ws.Range("E35").Value = '  -0.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '24.01'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.14%  '
$ws.Range("E37").Value = '  +0.73%  '
$ws.Range("E38").Value = '  +2.55%  '
$ws.Range("E39").Value = '  +0.93%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '166.81'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.55%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0789'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.42%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '26.65'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +9.54%  '
$ws.Range("E43").Value = '  +1.92%  '
$ws.Range("E44").Value = '  +1.86%  '
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '42.30'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.40%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.71'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.12%  '
$ws.Range("D48").Value = '2.604.95'
$ws.Range("E48").Value = '  +5.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.16'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.21%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.04'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.24%  '
$ws.Range("E51").Value = '  -0.46%  '
